# Auto-update predictions and index for 2025-10-22
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the "Win %" column to stay a text string ("72%" etc.) rather than
# being auto-coerced into a numeric percentage by the COM Value setter.
$ws.Range("F2:F15").NumberFormat = "@"

# ---- Row 2 ----
$ws.Range("B2").Value = "Al-Hilal SFC ✓ - Al-Sadd SC: 3:1"
$ws.Range("C2").Value = 3.49
$ws.Range("D2").Value = "Al-Hilal SFC"
$ws.Range("E2").Value = 4.5
$ws.Range("F2").Value = "73%"
$ws.Range("G2").Value = "✓"
$ws.Range("H2").Value = 4
$ws.Range("I2").Value = $true

# ---- Row 3 ----
$ws.Range("B3").Value = "Universitario de Deportes ✓ - Ayacucho FC: 2:1"
$ws.Range("C3").Value = 1.4
$ws.Range("D3").Value = "Universitario de Deportes"
$ws.Range("E3").Value = 2.5
$ws.Range("F3").Value = "72%"
$ws.Range("G3").Value = "✓"
$ws.Range("H3").Value = 3
$ws.Range("I3").Value = $false

# ---- Row 4 ----
$ws.Range("B4").Value = "FC Zbrojovka Brno  - Slezsky FC Opava: 0:0"
$ws.Range("C4").Value = 1.59
$ws.Range("D4").Value = "FC Zbrojovka Brno"
$ws.Range("E4").Value = 2.5
$ws.Range("F4").Value = "72%"
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = $true

# ---- Row 5 ----
$ws.Range("B5").Value = "Arsenal FC ✓ - Atlético de Madrid: 4:0"
$ws.Range("C5").Value = 1.98
$ws.Range("D5").Value = "Arsenal FC"
$ws.Range("E5").Value = 2.5
$ws.Range("F5").Value = "71%"
$ws.Range("G5").Value = "✓"
$ws.Range("H5").Value = 4
$ws.Range("I5").Value = $false

# ---- Row 6 ----
$ws.Range("B6").Value = "Union Saint-Gilloise - Inter Milan ✓: 0:4"
$ws.Range("C6").Value = 1.76
$ws.Range("D6").Value = "Inter Milan"
$ws.Range("E6").Value = 2.5
$ws.Range("F6").Value = "70%"
$ws.Range("G6").Value = "✓"
$ws.Range("H6").Value = 4
$ws.Range("I6").Value = $false

# ---- Row 7 ----
$ws.Range("B7").Value = "Chengdu Rongcheng - Johor Darul Ta'zim ✓: 0:2"
$ws.Range("C7").Value = 1.4
$ws.Range("D7").Value = "Johor Darul Ta'zim"
$ws.Range("E7").Value = 2.5
$ws.Range("F7").Value = "70%"
$ws.Range("G7").Value = "✓"
$ws.Range("H7").Value = 2
$ws.Range("I7").Value = $true

# ---- Row 8 ----
$ws.Range("B8").Value = "Lernayin Artsakh Goris X - FC Bentonit Ijevan: 0:6"
$ws.Range("C8").Value = 3.15
$ws.Range("D8").Value = "Lernayin Artsakh Goris"
$ws.Range("E8").Value = 4.5
$ws.Range("F8").Value = "68%"
$ws.Range("G8").Value = "X"
$ws.Range("H8").Value = 6
$ws.Range("I8").Value = $false

# ---- Row 9 ----
$ws.Range("B9").Value = "Bayer 04 Leverkusen - Paris Saint-Germain ✓: 2:7"
$ws.Range("C9").Value = 3.04
$ws.Range("D9").Value = "Paris Saint-Germain"
$ws.Range("E9").Value = 4.5
$ws.Range("F9").Value = "67%"
$ws.Range("G9").Value = "✓"
$ws.Range("H9").Value = 9
$ws.Range("I9").Value = $false

# ---- Row 10 (new) ----
$ws.Range("A10").Value = "Tue Oct 21"
$ws.Range("B10").Value = "FC Barcelona ✓ - Olympiacos Piraeus: 6:1"
$ws.Range("C10").Value = 2.44
$ws.Range("D10").Value = "FC Barcelona"
$ws.Range("E10").Value = 3.5
$ws.Range("F10").Value = "67%"
$ws.Range("G10").Value = "✓"
$ws.Range("H10").Value = 7
$ws.Range("I10").Value = $false

# ---- Row 11 (new) ----
$ws.Range("A11").Value = "Tue Oct 21"
$ws.Range("B11").Value = "FC Flora Tallinn ✓ - JK Trans Narva: 4:0"
$ws.Range("C11").Value = 2.37
$ws.Range("D11").Value = "FC Flora Tallinn"
$ws.Range("E11").Value = 3.5
$ws.Range("F11").Value = "60%"
$ws.Range("G11").Value = "✓"
$ws.Range("H11").Value = 4
$ws.Range("I11").Value = $false

# ---- Row 12 (new) ----
$ws.Range("A12").Value = "Tue Oct 21"
$ws.Range("B12").Value = "CSKA Moscow ✓ - Akron Togliatti: 3:2"
$ws.Range("C12").Value = 2.86
$ws.Range("D12").Value = "CSKA Moscow"
$ws.Range("E12").Value = 3.5
$ws.Range("F12").Value = "59%"
$ws.Range("G12").Value = "✓"
$ws.Range("H12").Value = 5
$ws.Range("I12").Value = $false

# ---- Row 13 (new) ----
$ws.Range("A13").Value = "Tue Oct 21"
$ws.Range("B13").Value = "Club Always Ready ✓ - Blooming Santa Cruz: 4:2"
$ws.Range("C13").Value = 2.73
$ws.Range("D13").Value = "Club Always Ready"
$ws.Range("E13").Value = 3.5
$ws.Range("F13").Value = "58%"
$ws.Range("G13").Value = "✓"
$ws.Range("H13").Value = 6
$ws.Range("I13").Value = $false

# ---- Row 14 (new) ----
$ws.Range("A14").Value = "Tue Oct 21"
$ws.Range("B14").Value = "MC Algiers  - Paradou AC: 23:00"
$ws.Range("C14").Value = 1.53
$ws.Range("D14").Value = "MC Algiers"
$ws.Range("E14").Value = 2.5
$ws.Range("F14").Value = "53%"
$ws.Range("H14").Value = 23
$ws.Range("I14").Value = $false

# ---- Row 15 (new) ----
$ws.Range("A15").Value = "Tue Oct 21"
$ws.Range("B15").Value = "Torpedo Kutaisi ✓ - Iberia 1999 Tbilisi: 3:1"
$ws.Range("C15").Value = 1.69
$ws.Range("D15").Value = "Torpedo Kutaisi"
$ws.Range("E15").Value = 2.5
$ws.Range("F15").Value = "52%"
$ws.Range("G15").Value = "✓"
$ws.Range("H15").Value = 4
$ws.Range("I15").Value = $false

# ---- Summary/index block moves from rows 10-12 to rows 16-18 ----
$ws.Range("K10").ClearContents()
$ws.Range("L10").ClearContents()
$ws.Range("K11").ClearContents()
$ws.Range("K12").ClearContents()

$ws.Range("K16").Formula = "=COUNTIF(I:I,TRUE)"
$ws.Range("L16").Formula = "=(K16/K18)*100"
$ws.Range("K17").Formula = "=COUNTIF(I:I,FALSE)"
$ws.Range("K18").Formula = "=K16+K17"
